$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.076.60"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.37%  '

$ws.Range("D3").Value = "'1.799.72"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.06%  '

$ws.Range("D4").Value = "'1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").Value = "'309.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.21%  '

$ws.Range("D6").Value = "'1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.09%  '

$ws.Range("D7").Value = "'0.5090"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.47%  '

$ws.Range("D8").Value = "'0.3842"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.77%  '

$ws.Range("D9").Value = "'0.07694"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.50%  '

$ws.Range("D10").Value = "'1.096"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.02%  '

$ws.Range("D11").Value = "'40.76"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.50%  '

$ws.Range("D12").Value = "'6.332"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.10%  '

$ws.Range("D13").Value = "'1.003"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.02%  '

$ws.Range("D14").Value = "'20.29"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.59%  '

$ws.Range("D15").Value = "'1.801.13"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.30%  '

$ws.Range("D16").Value = "'7.256"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.87%  '

$ws.Range("D17").Value = "'92.10"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.60%  '

$ws.Range("D18").Value = "'0.00001069"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.34%  '

$ws.Range("D19").Value = "'0.06560"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.85%  '

$ws.Range("E20").Value = '  -0.11%  '

$ws.Range("D21").Value = "'17.22"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.74%  '

$ws.Range("D22").Value = "'5.957"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.31%  '

$ws.Range("D23").Value = "'28.107.38"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.42%  '

$ws.Range("D24").Value = "'11.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.27%  '

$ws.Range("D25").Value = "'2.239"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.03%  '

$ws.Range("D26").Value = "'159.75"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.52%  '

$ws.Range("B27").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C27").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D27").Value = "'2.011.81"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.16%  '

$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D28").Value = "'2.409"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.56%  '

$ws.Range("D29").Value = "'20.21"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.58%  '

$ws.Range("D30").Value = "'127.40"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.65%  '

$ws.Range("D31").Value = "'0.1085"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.82%  '

$ws.Range("D32").Value = "'1.043"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.52%  '

$ws.Range("D33").Value = "'3.652"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.30%  '

$ws.Range("D34").Value = "'5.531"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.09%  '

$ws.Range("D35").Value = "'0.06959"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.20%  '

$ws.Range("D36").Value = "'9.004"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.10%  '

$ws.Range("D37").Value = "'0.02333"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.97%  '

$ws.Range("D38").Value = "'0.2163"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.21%  '

$ws.Range("D39").Value = "'4.993"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.21%  '

$ws.Range("D40").Value = "'11.42"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.43%  '

$ws.Range("D41").Value = "'0.6103"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.38%  '

$ws.Range("E42").Value = '  -0.13%  '

$ws.Range("E43").Value = '  -1.53%  '

$ws.Range("D44").Value = "'13.21"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.61%  '

$ws.Range("D45").Value = "'1.296"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.33%  '

$ws.Range("D46").Value = "'0.5889"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.58%  '

$ws.Range("D47").Value = "'3.714"
$ws.Range("D47").Style = "Normal"

$ws.Range("D48").Value = "'125.24"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.05%  '

$ws.Range("D49").Value = "'1.184"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.34%  '

$ws.Range("D50").Value = "'1.918"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.44%  '

$ws.Range("D51").Value = "'0.06733"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.32%  '
